# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet -- this pushes the old "Late" / heading / "Outstanding" columns
# one slot to the right (N->O, O->P, P->Q) and leaves the freshly
# inserted column N empty, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at N; existing N:P data (and the header labels)
# shift right to O:Q automatically.
$ws.Columns("N").Insert() | Out-Null

# The newly inserted column picks up the default sheet width; give it the
# same explicit width the author's saved file shows (~11 chars).
$ws.Columns("N").ColumnWidth = 10.15

# Bring this sheet to the front and leave the cursor where the author
# left it before saving.
$ws.Activate() | Out-Null
$ws.Range("S9").Select() | Out-Null
